# Update test from AI
# Applies the TestData.xlsx edits: new H-column weights on rows 7-11,
# shrunk row heights for rows 7-8, four new "question type" rows (18-21)
# styled in Arial/#111827, refreshed selection/dimension, and a portrait
# page setup on Sheet2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# --- Row height tweaks (row 7: 135 -> 75, row 8: 90 -> 60) ---
$ws.Rows.Item(7).RowHeight = 75
$ws.Rows.Item(8).RowHeight = 60

# --- New "weight" values in column H for rows 7-11 ---
$ws.Range("H7").Value = 5
$ws.Range("H8").Value = 5
$ws.Range("H9").Value = 5
$ws.Range("H10").Value = 5
$ws.Range("H11").Value = 6

# --- Four new rows listing question-type labels ---
$ws.Range("A18").Value = "Essay-SQA-Lead"
$ws.Range("A19").Value = "Case-SQA-Lead"
$ws.Range("A20").Value = "wert"
$ws.Range("A21").Value = "test-code-004"

# Give the new rows an Arial / #111827 font via a throwaway named style
# (keeps the cellStyles/cellStyleXfs tables untouched while still landing
# a fresh cellXfs entry that applies the new font).
$fontStyle = $wb.Styles.Add("TempFontStyle")
$fontStyle.Font.Name = "Arial"
$fontStyle.Font.Color = 2562065
$ws.Range("A18:A21").Style = "TempFontStyle"
$wb.Styles.Item("TempFontStyle").Delete()

# --- Selection / scroll position ---
$ws.Range("C27").Select()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

$wb.Save()
